$wb = $excel.ActiveWorkbook

# Add "Uygulama" sheet right after "Personel"
$personel = $wb.Worksheets.Item("Personel")
$sheetUygulama = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $personel)
$sheetUygulama.Name = "Uygulama"
$sheetUygulama.Range("A1").Value = "GÜNCELLEME"
$sheetUygulama.Range("A2").Value = "2/B"

# Add "İlçe" sheet right after "Uygulama" - this becomes the active sheet/tab
$sheetIlce = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheetUygulama)
$sheetIlce.Name = "İlçe"

Write-Output "Sheets:"
foreach ($s in $wb.Worksheets) {
    Write-Output $s.Name
}
